$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the whole used range first (the sheet shrinks from a 6x6 matrix
# to a 4x4 matrix: demand2 / net2 rows+cols are dropped).
$ws.Cells.Clear()

# Write cell-by-cell in strict row-major (top-to-bottom, left-to-right)
# order so new shared-string entries land in the same order Excel would
# create them while scanning/saving the sheet.
$ws.Cells.Item(1, 2).Value = "P_from_demand1"
$ws.Cells.Item(1, 3).Value = "P_from_net1"
$ws.Cells.Item(1, 4).Value = "P_from_pv1"
$ws.Cells.Item(1, 5).Value = "P_from_bat1"

$ws.Cells.Item(2, 1).Value = "P_to_demand1"
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = "P_net1_demand1"
$ws.Cells.Item(2, 4).Value = "P_pv1_demand1"
$ws.Cells.Item(2, 5).Value = "P_bat1_demand1"

$ws.Cells.Item(3, 1).Value = "P_to_net1"
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = "P_pv1_net1"
$ws.Cells.Item(3, 5).Value = "P_bat1_net1"

$ws.Cells.Item(4, 1).Value = "P_to_pv1"
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0

$ws.Cells.Item(5, 1).Value = "P_to_bat1"
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = "P_net1_bat1"
$ws.Cells.Item(5, 4).Value = "P_pv1_bat1"
$ws.Cells.Item(5, 5).Value = 0

# Bold / centered / bordered style for header row + row-label column,
# matching the original table's header styling.
$headerStyle = $ws.Range("B1:E1")
$headerStyle.Font.Bold = $true
$headerStyle.HorizontalAlignment = -4108
$headerStyle.VerticalAlignment = -4160
$headerStyle.Borders.LineStyle = 1

$labelStyle = $ws.Range("A2:A5")
$labelStyle.Font.Bold = $true
$labelStyle.HorizontalAlignment = -4108
$labelStyle.VerticalAlignment = -4160
$labelStyle.Borders.LineStyle = 1
